$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A63").Value = 62
$ws.Range("B63").Value = "30/05/2020"
$ws.Range("C63").Value = 137
$ws.Range("D63").Value = 4
$ws.Range("E63").Value = 32
$ws.Range("F63").Value = "90,19804066"
$ws.Range("G63").Value = "0,02919708029"
$ws.Range("H63").Value = 101
$ws.Range("I63").Value = 285
$ws.Range("J63").Value = 422
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 36
$ws.Range("M63").Value = 3
$ws.Range("N63").Value = 33
$ws.Range("O63").Value = 1
$ws.Range("P63").Value = 9

$ws.Range("A64").Value = 63
$ws.Range("B64").Value = "31/05/2020"
$ws.Range("C64").Value = 138
$ws.Range("D64").Value = 4
$ws.Range("E64").Value = 24
$ws.Range("F64").Value = "90,85642052"
$ws.Range("G64").Value = "0,02898550725"
$ws.Range("H64").Value = 110
$ws.Range("I64").Value = 287
$ws.Range("J64").Value = 425
$ws.Range("K64").Value = 1
$ws.Range("L64").Value = 35
$ws.Range("M64").Value = 2
$ws.Range("N64").Value = 33
$ws.Range("O64").Value = 3
$ws.Range("P64").Value = 10
